$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Home rails (final)" section entirely - from the
#    "Home rails (final)" heading paragraph through the end of the
#    "No Pro/Free split or outline here ..." paragraph (inclusive),
#    leaving the trailing empty paragraph untouched.
# ------------------------------------------------------------------

$startRange = $d.Content.Duplicate
$null = $startRange.Find.Execute("Home rails (final)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $startRange.Expand(4)

$endRange = $d.Content.Duplicate
$null = $endRange.Find.Execute("No Pro/Free split or outline here", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $endRange.Expand(4)

$sectionRange = $d.Range($startRange.Start, $endRange.End)
$sectionRange.Delete()

# ------------------------------------------------------------------
# 2) Mark the built-in "Default Paragraph Font" character style as
#    semi-hidden (adds <w:semiHidden/> in styles.xml). The Style
#    object in this automation surface does not expose a settable
#    SemiHidden property, so this is attempted defensively and will
#    not abort the rest of the script if unsupported.
# ------------------------------------------------------------------

$dpf = $d.Styles.Item("Default Paragraph Font")
try {
    $dpf.SemiHidden = $true
} catch {
    Write-Host "Style.SemiHidden not settable via this automation surface: $_"
}
